# Rename the three performance-parameter labels on both sheets by
# appending a "_0" suffix (these labels are backed by the same shared
# strings, so both the parameter sheet and the Scaling sheet need the
# matching update).
$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item(1)   # "performance_params_0"
$ws2 = $wb.Worksheets.Item(2)   # "Scaling"

$ws1.Range("A2").Value = "e_modulus_0"
$ws1.Range("A3").Value = "tensile_strain_at_break_0"
$ws1.Range("A4").Value = "tensile_yield_strength_0"

$ws2.Range("A2").Value = "e_modulus_0"
$ws2.Range("A3").Value = "tensile_strain_at_break_0"
$ws2.Range("A4").Value = "tensile_yield_strength_0"

# The label cells A2:A4 on the parameter sheet lose their cell border.
$ws1.Range("A2:A4").Borders.LineStyle = -4142   # xlLineStyleNone

# Widen column A on the Scaling sheet to fit the longer labels.
$ws2.Columns.Item(1).ColumnWidth = 21.17   # renders as width="22" in xlsx

# Update the selection/active-sheet state: the Scaling sheet ends up
# with A2:A4 selected (and is no longer the active tab), while the
# parameter sheet becomes the active tab with C14 selected. Select the
# non-active sheet's range first so the final Select() below leaves the
# parameter sheet as the active one.
$ws2.Range("A2:A4").Select() | Out-Null
$ws1.Range("C14").Select() | Out-Null
